$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.265.68"
$ws.Range("E2").Value = "  -5.61%  "
# Row 3
$ws.Range("D3").Value = "2.452.63"
$ws.Range("E3").Value = "  -8.20%  "
# Row 4
$ws.Range("E4").Value = "  -0.07%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.96"
$ws.Range("E5").Value = "  -3.19%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.42"
$ws.Range("E6").Value = "  -6.54%  "
# Row 7
$ws.Range("E7").Value = "  -0.33%  "
# Row 8
$ws.Range("E8").Value = "  -3.87%  "
# Row 9
$ws.Range("D9").Value = "2.466.84"
$ws.Range("E9").Value = "  -7.56%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0993"
$ws.Range("E10").Value = "  -5.80%  "
# Row 11
$ws.Range("E11").Value = "  -2.34%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.35"
$ws.Range("E12").Value = "  +1.42%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.350"
$ws.Range("E13").Value = "  -4.49%  "
# Row 14
$ws.Range("D14").Value = "2.882.90"
$ws.Range("E14").Value = "  -8.28%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.96"
$ws.Range("E15").Value = "  -8.11%  "
# Row 16
$ws.Range("D16").Value = "59.190.12"
$ws.Range("E16").Value = "  -5.51%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("E17").Value = "  -5.91%  "
# Row 18
$ws.Range("D18").Value = "2.506.86"
$ws.Range("E18").Value = "  -6.21%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.10"
$ws.Range("E19").Value = "  -6.33%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.34"
$ws.Range("E20").Value = "  -5.40%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.09"
$ws.Range("E21").Value = "  -6.12%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.967"
$ws.Range("E22").Value = "  -3.05%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.74"
$ws.Range("E23").Value = "  -8.31%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.461"
$ws.Range("E24").Value = "  -8.62%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "60.39"
$ws.Range("E25").Value = "  -4.04%  "
# Row 26
$ws.Range("E26").Value = "  -4.06%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.977"
$ws.Range("E27").Value = "  -2.31%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("E28").Value = "  -5.37%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.79"
$ws.Range("E29").Value = "  -6.09%  "
# Row 30
$ws.Range("E30").Value = "  -6.19%  "
# Row 31
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0770"
$ws.Range("E31").Value = "  -10.16%  "
# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.25"
$ws.Range("E32").Value = "  -8.28%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.997"
$ws.Range("E33").Value = "  -0.13%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.46"
$ws.Range("E34").Value = "  -6.54%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.40"
$ws.Range("E35").Value = "  -5.87%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("E36").Value = "  -6.31%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.36"
$ws.Range("E37").Value = "  -5.71%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.73"
$ws.Range("E38").Value = "  -2.47%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.83"
$ws.Range("E39").Value = "  -6.72%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "312.02"
$ws.Range("E40").Value = "  -9.78%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.65"
$ws.Range("E41").Value = "  -4.10%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.833"
$ws.Range("E42").Value = "  -12.87%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.71"
$ws.Range("E43").Value = "  -7.00%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.35%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.72"
$ws.Range("E45").Value = "  -2.77%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.583"
$ws.Range("E46").Value = "  -5.23%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0938"
$ws.Range("E47").Value = "  -3.36%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0524"
$ws.Range("E48").Value = "  -6.80%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0228"
$ws.Range("E49").Value = "  -5.31%  "
# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.47"
$ws.Range("E50").Value = "  -9.11%  "
# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.70"
$ws.Range("E51").Value = "  -5.51%  "
